$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.160.85"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.847.78"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "361.25"
$ws.Range("E5").Value = "  +5.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.63"
$ws.Range("E6").Value = "  -3.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  +3.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  +3.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.58"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.01"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "3.298.56"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "2.856.24"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.905"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").Value = "51.998.30"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  +8.91%  "
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").Value = "0.0₃0993"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.34"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.62"
$ws.Range("E24").Value = "  -3.89%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.42"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "53.71"
$ws.Range("E30").Value = "  +6.54%  "
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0465"
$ws.Range("E32").Value = "  +24.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.17"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.91"
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.40"
$ws.Range("E35").Value = "  +8.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0844"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.37"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.89"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "128.21"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.56"
$ws.Range("E44").Value = "  -7.93%  "
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("D46").Value = "2.115.35"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +8.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.86"
$ws.Range("E50").Value = "  +5.54%  "
$ws.Range("E51").Value = "  +1.30%  "
